# Update quiz results sheet:
#  - Row 2: ID changes from 37 to 79 (rest of row unchanged)
#  - Row 3: becomes a "Deepa" retake entry (ID 80, new duration/date)
#  - Row 4 (new): "Ajay" entry that used to be in row 3 (ID 81)
#  - Row 5 (new): brand-new "Test User" entry (ID 82)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - only the ID (column A) changes
$ws.Cells.Item(2, 1).Value = 79

# Row 3 - becomes Deepa's retake result
$ws.Cells.Item(3, 1).Value = 80
$ws.Cells.Item(3, 2).Value = "Deepa"
$ws.Cells.Item(3, 3).Value = "I079692"
$ws.Cells.Item(3, 4).Value = 100
$ws.Cells.Item(3, 5).Value = 5
$ws.Cells.Item(3, 6).Value = 5
$ws.Cells.Item(3, 7).Value = 21
$ws.Cells.Item(3, 8).Value = "2025-04-27 15:47:04"

# Row 4 (new) - Ajay's entry (previously in row 3)
$ws.Cells.Item(4, 1).Value = 81
$ws.Cells.Item(4, 2).Value = "Ajay"
$ws.Cells.Item(4, 3).Value = "I05235"
$ws.Cells.Item(4, 4).Value = 100
$ws.Cells.Item(4, 5).Value = 5
$ws.Cells.Item(4, 6).Value = 5
$ws.Cells.Item(4, 7).Value = 26
$ws.Cells.Item(4, 8).Value = "2025-04-27 13:06:43"

# Row 5 (new) - Test User entry
$ws.Cells.Item(5, 1).Value = 82
$ws.Cells.Item(5, 2).Value = "Test User"
$ws.Cells.Item(5, 3).Value = "I999999"
$ws.Cells.Item(5, 4).Value = 80
$ws.Cells.Item(5, 5).Value = 5
$ws.Cells.Item(5, 6).Value = 4
$ws.Cells.Item(5, 7).Value = 60
$ws.Cells.Item(5, 8).Value = "2025-04-27 15:43:16"
